$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.262727499008179
$ws.Range("B1").Value = 2.513934850692749
$ws.Range("C1").Value = 3.590359210968018
$ws.Range("D1").Value = 2.966970443725586
$ws.Range("E1").Value = 1.071547508239746
